$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its values as text, since many
# prices look like numbers/dates (e.g. "0.9982", "30.110.25") and Excel
# would otherwise auto-convert them.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.110.25'
$ws.Range("E2").Value = '  +8.33%  '
$ws.Range("D3").Value = '1.872.99'
$ws.Range("E3").Value = '  +5.47%  '
$ws.Range("D4").Value = '0.9982'
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = '249.66'
$ws.Range("E5").Value = '  +2.70%  '
$ws.Range("D6").Value = '0.9981'
$ws.Range("E6").Value = '  +0.22%  '
$ws.Range("D7").Value = '0.4976'
$ws.Range("E7").Value = '  +2.53%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").Value = '0.2847'
$ws.Range("E8").Value = '  +7.42%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").Value = '0.06586'
$ws.Range("E9").Value = '  +5.83%  '
$ws.Range("B10").Value = 'WrappedEther'
$ws.Range("C10").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D10").Value = '1.869.60'
$ws.Range("E10").Value = '  +5.29%  '
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").Value = '17.00'
$ws.Range("E11").Value = '  +4.53%  '
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").Value = '0.07200'
$ws.Range("E12").Value = '  +3.43%  '
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '0.6614'
$ws.Range("E13").Value = '  +7.65%  '
$ws.Range("B14").Value = 'Litecoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D14").Value = '85.02'
$ws.Range("E14").Value = '  +8.29%  '
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").Value = '4.808'
$ws.Range("E15").Value = '  +5.35%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '30.059.87'
$ws.Range("E16").Value = '  +8.35%  '
$ws.Range("B17").Value = 'Dai'
$ws.Range("C17").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D17").Value = '1.000'
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").Value = '12.86'
$ws.Range("E18").Value = '  +10.23%  '
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").Value = '0.000007510'
$ws.Range("E19").Value = '  +5.41%  '
$ws.Range("B20").Value = 'BinanceUSD'
$ws.Range("C20").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D20").Value = '0.9982'
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").Value = '2.108.07'
$ws.Range("E21").Value = '  +6.30%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").Value = '4.752'
$ws.Range("E22").Value = '  +4.14%  '
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").Value = '9.047'
$ws.Range("E23").Value = '  +5.90%  '
$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").Value = '5.502'
$ws.Range("E24").Value = '  +6.26%  '
$ws.Range("B25").Value = 'Monero'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D25").Value = '145.14'
$ws.Range("E25").Value = '  +2.63%  '
$ws.Range("B26").Value = 'BitcoinCash'
$ws.Range("C26").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D26").Value = '136.72'
$ws.Range("E26").Value = '  +24.03%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '16.74'
$ws.Range("E27").Value = '  +7.84%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = '1.947'
$ws.Range("E28").Value = '  +3.77%  '
$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").Value = '1.390'
$ws.Range("E29").Value = '  +0.38%  '
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = '4.245'
$ws.Range("E30").Value = '  +5.97%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").Value = '0.08632'
$ws.Range("E31").Value = '  +4.89%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").Value = '3.898'
$ws.Range("E32").Value = '  +4.47%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").Value = '0.05077'
$ws.Range("E33").Value = '  +8.43%  '
$ws.Range("B34").Value = 'ARBITRUM'
$ws.Range("C34").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D34").Value = '1.130'
$ws.Range("E34").Value = '  +8.30%  '
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '0.6857'
$ws.Range("E35").Value = '  +8.90%  '
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D36").Value = '2.689'
$ws.Range("E36").Value = '  +2.84%  '
$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").Value = '2.317'
$ws.Range("E37").Value = '  +13.62%  '
$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").Value = '2.738'
$ws.Range("E38").Value = '  +5.31%  '
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").Value = '0.9617'
$ws.Range("E39").Value = '  +3.06%  '
$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").Value = '0.01629'
$ws.Range("E40").Value = '  +7.22%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '6.115'
$ws.Range("E41").Value = '  +6.32%  '
$ws.Range("B42").Value = 'PaxDollar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D42").Value = '0.9990'
$ws.Range("E42").Value = '  +0.25%  '
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").Value = '103.76'
$ws.Range("E43").Value = '  +4.01%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.4203'
$ws.Range("E44").Value = '  +7.56%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '7.429'
$ws.Range("E45").Value = '  +6.81%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '0.1254'
$ws.Range("E46").Value = '  +6.75%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '0.05622'
$ws.Range("E47").Value = '  +4.26%  '
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").Value = '32.54'
$ws.Range("E48").Value = '  +7.68%  '
$ws.Range("D49").Value = '8.321'
$ws.Range("E49").Value = '  +5.54%  '
$ws.Range("B50").Value = 'Decentraland'
$ws.Range("C50").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D50").Value = '0.3732'
$ws.Range("E50").Value = '  +8.67%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.342'
$ws.Range("E51").Value = '  +5.94%  '
